$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: current_phase 1 -> 2
$ws.Range("D29").Value = 2

# Row 30: current_phase 1 -> 2
$ws.Range("D30").Value = 2

# Row 30: last_action_date set (stored as text string)
$ws.Range("E30").Value = "2026-02-12T19:50:44.100273+00:00"

# Row 30: replies_count 0 -> 1
$ws.Range("I30").Value = 1

# Row 30: replied_message_ids [] -> [10]
$ws.Range("M30").Value = "[10]"
